$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in the title cell
$ws.Range("A1").Value = "Datos actualizados a 23 de Junio de 2020 a las 21:07"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 2409068
$ws.Range("C4").Value = 20915
$ws.Range("D4").Value = 1005688
$ws.Range("E4").Value = 1280167
$ws.Range("G4").Value = 603
$ws.Range("H4").Value = 123213

# Row 7: India -> India
$ws.Range("B7").Value = 455859
$ws.Range("C7").Value = 15409
$ws.Range("E7").Value = 182853

# Row 14: Alemania -> Alemania
$ws.Range("B14").Value = 192539
$ws.Range("C14").Value = 420
$ws.Range("E14").Value = 7860

# Row 19: Francia -> Francia
$ws.Range("B19").Value = 161267
$ws.Range("C19").Value = 517
$ws.Range("D19").Value = 74871
$ws.Range("E19").Value = 56676
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = 29720

# Row 30: Ecuador -> Ecuador
$ws.Range("E30").Value = 21375
$ws.Range("G30").Value = 51
$ws.Range("H30").Value = 4274

# Row 50: Barein -> Barein
$ws.Range("E50").Value = 5478
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 67

# Row 76: Uzbekistan -> Uzbekistan
$ws.Range("B76").Value = 6662
$ws.Range("C76").Value = 201
$ws.Range("E76").Value = 2083

# Row 95: Tailandia -> Mauritania
$ws.Range("A95").Value = "Mauritania"
$ws.Range("B95").Value = 3292
$ws.Range("C95").Value = 171
$ws.Range("D95").Value = 963
$ws.Range("E95").Value = 2215
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 114

# Row 96: Mauritania -> Tailandia
$ws.Range("A96").Value = "Tailandia"
$ws.Range("B96").Value = 3156
$ws.Range("C96").Value = 5
$ws.Range("D96").Value = 3023
$ws.Range("E96").Value = 75
$ws.Range("H96").Value = 58

# Row 124: Tunez -> Estado de Palestina
$ws.Range("A124").Value = "Estado de Palestina"
$ws.Range("B124").Value = 1169
$ws.Range("C124").Value = 168
$ws.Range("D124").Value = 442
$ws.Range("E124").Value = 724
$ws.Range("H124").Value = 3

# Row 125: Estado de Palestina -> Tunez
$ws.Range("A125").Value = "Tunez"
$ws.Range("B125").Value = 1159
$ws.Range("C125").Value = 0
$ws.Range("D125").Value = 1023
$ws.Range("E125").Value = 86
$ws.Range("H125").Value = 50

# Row 130: Republica de Chipre -> Yemen
$ws.Range("A130").Value = "Yemen"
$ws.Range("B130").Value = 992
$ws.Range("C130").Value = 25
$ws.Range("D130").Value = 356
$ws.Range("E130").Value = 375
$ws.Range("G130").Value = 4
$ws.Range("H130").Value = 261

# Row 131: Yemen -> Republica de Chipre
$ws.Range("A131").Value = "Republica de Chipre"
$ws.Range("B131").Value = 990
$ws.Range("C131").Value = 2
$ws.Range("D131").Value = 824
$ws.Range("E131").Value = 147
$ws.Range("H131").Value = 19

# Row 133: Georgia -> Georgia
$ws.Range("D133").Value = 768
$ws.Range("E133").Value = 129

# Row 136: Republica del Chad -> Republica del Chad
$ws.Range("B136").Value = 860
$ws.Range("C136").Value = 2
$ws.Range("D136").Value = 757

# Row 146: Jamaica -> Suazilandia
$ws.Range("A146").Value = "Suazilandia"
$ws.Range("B146").Value = 674
$ws.Range("C146").Value = 31
$ws.Range("D146").Value = 319
$ws.Range("E146").Value = 348
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 7

# Row 147: Malta -> Jamaica
$ws.Range("A147").Value = "Jamaica"
$ws.Range("C147").Value = 6
$ws.Range("D147").Value = 516
$ws.Range("E147").Value = 139
$ws.Range("H147").Value = 10

# Row 148: Liberia -> Malta
$ws.Range("A148").Value = "Malta"
$ws.Range("B148").Value = 665
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 618
$ws.Range("E148").Value = 38
$ws.Range("H148").Value = 9

# Row 149: Suazilandia -> Liberia
$ws.Range("A149").Value = "Liberia"
$ws.Range("B149").Value = 652
$ws.Range("C149").Value = 2
$ws.Range("D149").Value = 270
$ws.Range("E149").Value = 348
$ws.Range("H149").Value = 34

# Row 164: Siria -> Siria
$ws.Range("D164").Value = 94
$ws.Range("E164").Value = 118

# Row 168: Islas Feroe -> Angola
$ws.Range("A168").Value = "Angola"
$ws.Range("B168").Value = 189
$ws.Range("C168").Value = 3
$ws.Range("D168").Value = 77
$ws.Range("E168").Value = 102
$ws.Range("H168").Value = 10

# Row 169: Angola -> Islas Feroe
$ws.Range("A169").Value = "Islas Feroe"
$ws.Range("B169").Value = 187
$ws.Range("D169").Value = 187
$ws.Range("E169").Value = 0
$ws.Range("H169").Value = 0

# Row 185: Butan -> Namibia
$ws.Range("A185").Value = "Namibia"
$ws.Range("B185").Value = 72
$ws.Range("C185").Value = 9
$ws.Range("D185").Value = 21
$ws.Range("E185").Value = 51

# Row 186: Namibia -> Butan
$ws.Range("A186").Value = "Butan"
$ws.Range("B186").Value = 69
$ws.Range("C186").Value = 1
$ws.Range("D186").Value = 32
$ws.Range("E186").Value = 37
